$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.813634
$ws.Range("H2").Value = 2.440902
$ws.Range("I2").Value = 0.4058779337539379
$ws.Range("J2").Value = 0.405877933753938
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.05883533333333333
$ws.Range("N2").Value = 0.176506
$ws.Range("Q2").Value = 0.04787042760133333
$ws.Range("R2").Value = 0.430833848412
$ws.Range("S2").Value = 0.4058779337539379
$ws.Range("T2").Value = 0.405877933753938

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.009547
$ws.Range("H3").Value = 3.028641
$ws.Range("I3").Value = 0.5036083182210759
$ws.Range("J3").Value = 0.5036083182210759
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.05883533333333333
$ws.Range("N3").Value = 0.176506
$ws.Range("Q3").Value = 0.05939703426066668
$ws.Range("R3").Value = 0.534573308346
$ws.Range("S3").Value = 0.5036083182210759
$ws.Range("T3").Value = 0.5036083182210759

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.067537
$ws.Range("H4").Value = 0.202611
$ws.Range("I4").Value = 0.03369055129448831
$ws.Range("J4").Value = 0.03369055129448831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.05883533333333333
$ws.Range("N4").Value = 0.176506
$ws.Range("Q4").Value = 0.003973561907333333
$ws.Range("R4").Value = 0.035762057166
$ws.Range("S4").Value = 0.03369055129448831
$ws.Range("T4").Value = 0.03369055129448831

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1139093333333333
$ws.Range("H5").Value = 0.341728
$ws.Range("I5").Value = 0.05682319673049787
$ws.Range("J5").Value = 0.05682319673049788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.05883533333333333
$ws.Range("N5").Value = 0.176506
$ws.Range("Q5").Value = 0.006701893596444445
$ws.Range("R5").Value = 0.06031704236800001
$ws.Range("S5").Value = 0.05682319673049787
$ws.Range("T5").Value = 0.05682319673049788
